$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 826, shifting existing rows 826:867 down to 827:868
$ws.Rows.Item(826).Insert()

# Populate the newly inserted row 826 with the new data point.
# Force column A to be treated as plain text so the date-like string is not
# auto-converted into a date serial number, then restore the plain/default
# style so the cell matches the formatting of its neighbours.
$ws.Cells.Item(826, 1).NumberFormat = "@"
$ws.Cells.Item(826, 1).Value = "2026/02/17"
$ws.Cells.Item(826, 1).Style = "Normal"
$ws.Cells.Item(826, 2).Value = "火"
$ws.Cells.Item(826, 3).Value = 15
$ws.Cells.Item(826, 4).Value = 201
